# Auto-generated Excel COM-interop script
# Applies updated market-price / profit values across multiple crafting-leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2982.1538
$ws.Cells.Item(28, 9).Value = 596.7368
$ws.Cells.Item(28, 10).Value = 9456.857
$ws.Cells.Item(28, 11).Value = 596.7368
$ws.Cells.Item(28, 12).Value = 9456.857
$ws.Cells.Item(28, 13).Value = -111.7368
$ws.Cells.Item(28, 14).Value = -10426.857

$ws.Cells.Item(74, 8).Value = 3942.0833
$ws.Cells.Item(74, 9).Value = 3845.9092
$ws.Cells.Item(74, 11).Value = 3845.9092
$ws.Cells.Item(74, 13).Value = -2909.9092

$ws.Cells.Item(77, 8).Value = 3942.0833
$ws.Cells.Item(77, 9).Value = 3845.9092
$ws.Cells.Item(77, 11).Value = 19229.546
$ws.Cells.Item(77, 13).Value = -14549.546

$ws.Cells.Item(80, 8).Value = 2587024
$ws.Cells.Item(80, 9).Value = 1226078
$ws.Cells.Item(80, 11).Value = 3678234
$ws.Cells.Item(80, 13).Value = -3677236

$ws.Cells.Item(83, 8).Value = 2587024
$ws.Cells.Item(83, 9).Value = 1226078
$ws.Cells.Item(83, 11).Value = 11034702
$ws.Cells.Item(83, 13).Value = -11029710

$ws.Cells.Item(92, 8).Value = 3023.6538
$ws.Cells.Item(92, 9).Value = 1169.7
$ws.Cells.Item(92, 10).Value = 9203.5
$ws.Cells.Item(92, 11).Value = 1169.7
$ws.Cells.Item(92, 12).Value = 9203.5
$ws.Cells.Item(92, 13).Value = 78.29999999999995
$ws.Cells.Item(92, 14).Value = -11699.5

$ws.Cells.Item(101, 8).Value = 296.33334
$ws.Cells.Item(101, 10).Value = 396
$ws.Cells.Item(101, 12).Value = 1188
$ws.Cells.Item(101, 14).Value = -4432

$ws.Cells.Item(106, 8).Value = 5585.9443
$ws.Cells.Item(106, 9).Value = 6874.846
$ws.Cells.Item(106, 11).Value = 6874.846
$ws.Cells.Item(106, 13).Value = -6243.846

$ws.Cells.Item(127, 8).Value = 1886
$ws.Cells.Item(127, 10).Value = 5000
$ws.Cells.Item(127, 12).Value = 15000
$ws.Cells.Item(127, 14).Value = -24920

$ws.Cells.Item(132, 8).Value = 1402.18
$ws.Cells.Item(132, 9).Value = 1015.95654
$ws.Cells.Item(132, 11).Value = 3047.86962
$ws.Cells.Item(132, 13).Value = -517.8696199999999

$ws.Cells.Item(135, 8).Value = 959
$ws.Cells.Item(135, 9).Value = 341.16666
$ws.Cells.Item(135, 11).Value = 3070.49994
$ws.Cells.Item(135, 13).Value = -535.4999399999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3541.5
$ws.Cells.Item(32, 9).Value = 3291.066
$ws.Cells.Item(32, 10).Value = 6073.6665
$ws.Cells.Item(32, 11).Value = 3291.066
$ws.Cells.Item(32, 12).Value = 6073.6665
$ws.Cells.Item(32, 13).Value = -3004.066
$ws.Cells.Item(32, 14).Value = -6647.6665

$ws.Cells.Item(45, 8).Value = 1289419.9
$ws.Cells.Item(45, 9).Value = 1906490.1
$ws.Cells.Item(45, 11).Value = 1906490.1
$ws.Cells.Item(45, 13).Value = -1906113.1

$ws.Cells.Item(60, 8).Value = 89693.55
$ws.Cells.Item(60, 9).Value = 89693.55
$ws.Cells.Item(60, 11).Value = 89693.55
$ws.Cells.Item(60, 13).Value = -88960.55

$ws.Cells.Item(74, 8).Value = 1599.2174
$ws.Cells.Item(74, 9).Value = 1536.6111
$ws.Cells.Item(74, 11).Value = 1536.6111
$ws.Cells.Item(74, 13).Value = -662.6111000000001

$ws.Cells.Item(77, 8).Value = 1599.2174
$ws.Cells.Item(77, 9).Value = 1536.6111
$ws.Cells.Item(77, 11).Value = 7683.0555
$ws.Cells.Item(77, 13).Value = -3315.0555

$ws.Cells.Item(97, 8).Value = 5323.7646
$ws.Cells.Item(97, 9).Value = 5342.1665
$ws.Cells.Item(97, 10).Value = 5279.6
$ws.Cells.Item(97, 11).Value = 5342.1665
$ws.Cells.Item(97, 12).Value = 5279.6
$ws.Cells.Item(97, 13).Value = -4846.1665
$ws.Cells.Item(97, 14).Value = -6271.6

$ws.Cells.Item(122, 8).Value = 5667.8945
$ws.Cells.Item(122, 9).Value = 3422
$ws.Cells.Item(122, 10).Value = 8163.3335
$ws.Cells.Item(122, 11).Value = 10266
$ws.Cells.Item(122, 12).Value = 24490.0005
$ws.Cells.Item(122, 13).Value = -7816
$ws.Cells.Item(122, 14).Value = -29390.0005

$ws.Cells.Item(132, 8).Value = 3127208
$ws.Cells.Item(132, 9).Value = 2195.2666
$ws.Cells.Item(132, 11).Value = 6585.7998
$ws.Cells.Item(132, 13).Value = -4055.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2340.72
$ws.Cells.Item(94, 9).Value = 2524.375
$ws.Cells.Item(94, 10).Value = 2014.2222
$ws.Cells.Item(94, 11).Value = 2524.375
$ws.Cells.Item(94, 12).Value = 2014.2222
$ws.Cells.Item(94, 13).Value = -2073.375
$ws.Cells.Item(94, 14).Value = -2916.2222

$ws.Cells.Item(132, 8).Value = 110000
$ws.Cells.Item(132, 10).Value = 110000
$ws.Cells.Item(132, 12).Value = 110000
$ws.Cells.Item(132, 14).Value = -120120

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 21281804
$ws.Cells.Item(31, 9).Value = 40004624
$ws.Cells.Item(31, 10).Value = 5873.409
$ws.Cells.Item(31, 11).Value = 40004624
$ws.Cells.Item(31, 12).Value = 5873.409
$ws.Cells.Item(31, 13).Value = -40004329
$ws.Cells.Item(31, 14).Value = -6463.409

$ws.Cells.Item(34, 8).Value = 21281804
$ws.Cells.Item(34, 9).Value = 40004624
$ws.Cells.Item(34, 10).Value = 5873.409
$ws.Cells.Item(34, 11).Value = 40004624
$ws.Cells.Item(34, 12).Value = 5873.409
$ws.Cells.Item(34, 13).Value = -40004422
$ws.Cells.Item(34, 14).Value = -6277.409

$ws.Cells.Item(132, 8).Value = 1259.9387
$ws.Cells.Item(132, 9).Value = 1259.9387
$ws.Cells.Item(132, 11).Value = 3779.8161
$ws.Cells.Item(132, 13).Value = -1249.8161

$ws.Cells.Item(134, 8).Value = 1096.1818
$ws.Cells.Item(134, 9).Value = 1126.2
$ws.Cells.Item(134, 11).Value = 3378.6
$ws.Cells.Item(134, 13).Value = -843.6000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 519.0714
$ws.Cells.Item(5, 9).Value = 395.7143
$ws.Cells.Item(5, 10).Value = 642.4286
$ws.Cells.Item(5, 11).Value = 1187.1429
$ws.Cells.Item(5, 12).Value = 1927.2858
$ws.Cells.Item(5, 13).Value = -1075.1429
$ws.Cells.Item(5, 14).Value = -2151.2858

$ws.Cells.Item(37, 8).Value = 92498.22
$ws.Cells.Item(37, 10).Value = 92498.22
$ws.Cells.Item(37, 12).Value = 277494.66
$ws.Cells.Item(37, 14).Value = -277718.66

$ws.Cells.Item(60, 8).Value = 4428.875
$ws.Cells.Item(60, 10).Value = 16865.5
$ws.Cells.Item(60, 12).Value = 50596.5
$ws.Cells.Item(60, 14).Value = -51098.5

$ws.Cells.Item(97, 8).Value = 400.25
$ws.Cells.Item(97, 9).Value = 403
$ws.Cells.Item(97, 10).Value = 399.33334
$ws.Cells.Item(97, 11).Value = 1209
$ws.Cells.Item(97, 12).Value = 1198.00002
$ws.Cells.Item(97, 13).Value = -713
$ws.Cells.Item(97, 14).Value = -2190.00002

$ws.Cells.Item(122, 8).Value = 10012.828
$ws.Cells.Item(122, 10).Value = 613.3
$ws.Cells.Item(122, 12).Value = 5519.7
$ws.Cells.Item(122, 14).Value = -10419.7

$ws.Cells.Item(135, 8).Value = 519.0714
$ws.Cells.Item(135, 9).Value = 395.7143
$ws.Cells.Item(135, 10).Value = 642.4286
$ws.Cells.Item(135, 11).Value = 3561.4287
$ws.Cells.Item(135, 12).Value = 5781.8574
$ws.Cells.Item(135, 13).Value = -1026.4287
$ws.Cells.Item(135, 14).Value = -10851.8574

$ws.Cells.Item(137, 8).Value = 7399.5
$ws.Cells.Item(137, 10).Value = 9799
$ws.Cells.Item(137, 12).Value = 29397
$ws.Cells.Item(137, 14).Value = -39597

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9089.111000000001
$ws.Cells.Item(70, 9).Value = 6312.4443
$ws.Cells.Item(70, 11).Value = 6312.4443
$ws.Cells.Item(70, 13).Value = -6042.4443

$ws.Cells.Item(73, 8).Value = 9089.111000000001
$ws.Cells.Item(73, 9).Value = 6312.4443
$ws.Cells.Item(73, 11).Value = 6312.4443
$ws.Cells.Item(73, 13).Value = -5376.4443

$ws.Cells.Item(132, 8).Value = 10003669
$ws.Cells.Item(132, 9).Value = 3955.5715
$ws.Cells.Item(132, 11).Value = 11866.7145
$ws.Cells.Item(132, 13).Value = -9336.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1538.3636
$ws.Cells.Item(22, 10).Value = 1676.6666
$ws.Cells.Item(22, 12).Value = 1676.6666
$ws.Cells.Item(22, 14).Value = -2266.6666

$ws.Cells.Item(27, 8).Value = 1538.3636
$ws.Cells.Item(27, 10).Value = 1676.6666
$ws.Cells.Item(27, 12).Value = 1676.6666
$ws.Cells.Item(27, 14).Value = -1890.6666

$ws.Cells.Item(68, 8).Value = 2454032.8
$ws.Cells.Item(68, 9).Value = 4168930.2
$ws.Cells.Item(68, 10).Value = 4179.143
$ws.Cells.Item(68, 11).Value = 4168930.2
$ws.Cells.Item(68, 12).Value = 4179.143
$ws.Cells.Item(68, 13).Value = -4168181.2
$ws.Cells.Item(68, 14).Value = -5677.143

$ws.Cells.Item(71, 8).Value = 2454032.8
$ws.Cells.Item(71, 9).Value = 4168930.2
$ws.Cells.Item(71, 10).Value = 4179.143
$ws.Cells.Item(71, 11).Value = 20844651
$ws.Cells.Item(71, 12).Value = 20895.715
$ws.Cells.Item(71, 13).Value = -20840907
$ws.Cells.Item(71, 14).Value = -28383.715

$ws.Cells.Item(132, 8).Value = 5505.7144
$ws.Cells.Item(132, 9).Value = 3222
$ws.Cells.Item(132, 10).Value = 6128.5454
$ws.Cells.Item(132, 11).Value = 9666
$ws.Cells.Item(132, 12).Value = 18385.6362
$ws.Cells.Item(132, 13).Value = -7136
$ws.Cells.Item(132, 14).Value = -23445.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 716376.4399999999
$ws.Cells.Item(132, 9).Value = 2262.7778
$ws.Cells.Item(132, 11).Value = 6788.3334
$ws.Cells.Item(132, 13).Value = -4258.3334

$ws.Cells.Item(136, 8).Value = 836596.2
$ws.Cells.Item(136, 9).Value = 3183.5
$ws.Cells.Item(136, 10).Value = 1253302.5
$ws.Cells.Item(136, 11).Value = 9550.5
$ws.Cells.Item(136, 12).Value = 3759907.5
$ws.Cells.Item(136, 13).Value = -7000.5
$ws.Cells.Item(136, 14).Value = -3765007.5
